$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Through 2022-07-24" to "Through 2022-07-25"
$ws.Name = "Through 2022-07-25"

# Update the header label in I1 (shared string "2022 (through 07-24)")
$ws.Range("I1").Value = "2022 (through 07-25)"

# Update August total (row 8) and grand Total row (row 14) for column I
$ws.Range("I8").Value = 139
$ws.Range("I14").Value = 945
